$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the formula for B17 (sbp_payment for 2011)
$ws.Range("B17").Formula = "=43.58+43.58+43.58"

# Update the selection / view state to match the saved workbook
$ws.Range("B17").Select()
